$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 340
$ws.Range('B340').Value = 5169591
$ws.Range('F340').Value = 'Hearts'
$ws.Range('G340').Value = 'Ross County'
$ws.Range('H340').Value = 2
$ws.Range('I340').Value = 1
$ws.Range('J340').Value = 'H'
$ws.Range('K340').Value = 1.666
$ws.Range('L340').Value = 3.75
$ws.Range('M340').Value = 5
$ws.Range('N340').Value = 1.5
$ws.Range('O340').Value = 4.2
$ws.Range('P340').Value = 6.5
$ws.Range('Q340').Value = -1
$ws.Range('R340').Value = 1.9
$ws.Range('S340').Value = 1.95
$ws.Range('T340').Value = 2.5
$ws.Range('U340').Value = 1.85
$ws.Range('V340').Value = 2
$ws.Range('W340').Value = 0.5
$ws.Range('X340').Value = -1
$ws.Range('Y340').Value = -1
$ws.Range('Z340').Value = 0
$ws.Range('AA340').Value = -0
$ws.Range('AB340').Value = 0.8500000000000001
$ws.Range('AC340').Value = -1

# Row 341
$ws.Range('B341').Value = 5169276
$ws.Range('F341').Value = 'St Johnstone'
$ws.Range('G341').Value = 'Hibernian'
$ws.Range('H341').Value = 0
$ws.Range('I341').Value = 1
$ws.Range('J341').Value = 'A'
$ws.Range('K341').Value = 2.625
$ws.Range('L341').Value = 3.3
$ws.Range('M341').Value = 2.625
$ws.Range('N341').Value = 3.25
$ws.Range('O341').Value = 3.1
$ws.Range('P341').Value = 2.3
$ws.Range('Q341').Value = 0.25
$ws.Range('R341').Value = 1.85
$ws.Range('S341').Value = 2
$ws.Range('T341').Value = 2
$ws.Range('U341').Value = 1.925
$ws.Range('V341').Value = 1.925
$ws.Range('W341').Value = -1
$ws.Range('X341').Value = -1
$ws.Range('Y341').Value = 1.3
$ws.Range('Z341').Value = -1
$ws.Range('AA341').Value = 1
$ws.Range('AB341').Value = -1
$ws.Range('AC341').Value = 0.925

# Row 391
$ws.Range('B391').Value = 5169301
$ws.Range('F391').Value = 'Livingston'
$ws.Range('G391').Value = 'Ross County'
$ws.Range('H391').Value = 0
$ws.Range('I391').Value = 1
$ws.Range('J391').Value = 'A'
$ws.Range('K391').Value = 1.909
$ws.Range('L391').Value = 3.2
$ws.Range('M391').Value = 4
$ws.Range('N391').Value = 1.909
$ws.Range('O391').Value = 3.3
$ws.Range('P391').Value = 4.333
$ws.Range('Q391').Value = -0.5
$ws.Range('R391').Value = 1.875
$ws.Range('S391').Value = 1.975
$ws.Range('T391').Value = 2.25
$ws.Range('U391').Value = 1.975
$ws.Range('V391').Value = 1.875
$ws.Range('W391').Value = -1
$ws.Range('X391').Value = -1
$ws.Range('Y391').Value = 3.333
$ws.Range('Z391').Value = -1
$ws.Range('AA391').Value = 0.9750000000000001
$ws.Range('AB391').Value = -1
$ws.Range('AC391').Value = 0.875

# Row 392
$ws.Range('B392').Value = 5169302
$ws.Range('F392').Value = 'Rangers'
$ws.Range('G392').Value = 'St Mirren'
$ws.Range('H392').Value = 4
$ws.Range('I392').Value = 0
$ws.Range('J392').Value = 'H'
$ws.Range('K392').Value = 1.25
$ws.Range('L392').Value = 5.5
$ws.Range('M392').Value = 11
$ws.Range('N392').Value = 1.2
$ws.Range('O392').Value = 6
$ws.Range('P392').Value = 15
$ws.Range('Q392').Value = -2
$ws.Range('R392').Value = 2.05
$ws.Range('S392').Value = 1.8
$ws.Range('T392').Value = 3
$ws.Range('U392').Value = 1.825
$ws.Range('V392').Value = 2.025
$ws.Range('W392').Value = 0.2
$ws.Range('X392').Value = -1
$ws.Range('Y392').Value = -1
$ws.Range('Z392').Value = 1.05
$ws.Range('AA392').Value = -1
$ws.Range('AB392').Value = 0.825
$ws.Range('AC392').Value = -1

# Row 406
$ws.Range('B406').Value = 5169617
$ws.Range('F406').Value = 'Kilmarnock'
$ws.Range('G406').Value = 'Ross County'
$ws.Range('H406').Value = 1
$ws.Range('I406').Value = 0
$ws.Range('J406').Value = 'H'
$ws.Range('K406').Value = 2.1
$ws.Range('L406').Value = 3.2
$ws.Range('M406').Value = 3.6
$ws.Range('N406').Value = 1.909
$ws.Range('O406').Value = 3.3
$ws.Range('P406').Value = 4.5
$ws.Range('Q406').Value = -0.5
$ws.Range('R406').Value = 1.925
$ws.Range('S406').Value = 1.925
$ws.Range('T406').Value = 2.25
$ws.Range('U406').Value = 2.05
$ws.Range('V406').Value = 1.8
$ws.Range('W406').Value = 0.909
$ws.Range('X406').Value = -1
$ws.Range('Y406').Value = -1
$ws.Range('Z406').Value = 0.925
$ws.Range('AA406').Value = -1
$ws.Range('AB406').Value = -1
$ws.Range('AC406').Value = 0.8

# Row 407
$ws.Range('B407').Value = 5169308
$ws.Range('F407').Value = 'St Mirren'
$ws.Range('G407').Value = 'Dundee Utd'
$ws.Range('H407').Value = 2
$ws.Range('I407').Value = 1
$ws.Range('J407').Value = 'H'
$ws.Range('K407').Value = 1.909
$ws.Range('L407').Value = 3.25
$ws.Range('M407').Value = 4.5
$ws.Range('N407').Value = 2.05
$ws.Range('O407').Value = 3.3
$ws.Range('P407').Value = 3.75
$ws.Range('Q407').Value = -0.5
$ws.Range('R407').Value = 2
$ws.Range('S407').Value = 1.85
$ws.Range('T407').Value = 2.25
$ws.Range('U407').Value = 1.925
$ws.Range('V407').Value = 1.925
$ws.Range('W407').Value = 1.05
$ws.Range('X407').Value = -1
$ws.Range('Y407').Value = -1
$ws.Range('Z407').Value = 1
$ws.Range('AA407').Value = -1
$ws.Range('AB407').Value = 0.925
$ws.Range('AC407').Value = -1

# Row 434
$ws.Range('B434').Value = 5169628
$ws.Range('F434').Value = 'Hearts'
$ws.Range('G434').Value = 'Kilmarnock'
$ws.Range('H434').Value = 3
$ws.Range('I434').Value = 1
$ws.Range('J434').Value = 'H'
$ws.Range('K434').Value = 1.7
$ws.Range('L434').Value = 3.6
$ws.Range('M434').Value = 5
$ws.Range('N434').Value = 1.666
$ws.Range('O434').Value = 3.8
$ws.Range('P434').Value = 5
$ws.Range('Q434').Value = -0.75
$ws.Range('R434').Value = 1.8
$ws.Range('S434').Value = 2.05
$ws.Range('T434').Value = 2.5
$ws.Range('U434').Value = 1.825
$ws.Range('V434').Value = 2.025
$ws.Range('W434').Value = 0.6659999999999999
$ws.Range('X434').Value = -1
$ws.Range('Y434').Value = -1
$ws.Range('Z434').Value = 0.8
$ws.Range('AA434').Value = -1
$ws.Range('AB434').Value = 0.825
$ws.Range('AC434').Value = -1

# Row 435
$ws.Range('B435').Value = 5169322
$ws.Range('F435').Value = 'Ross County'
$ws.Range('G435').Value = 'St Johnstone'
$ws.Range('H435').Value = 1
$ws.Range('I435').Value = 2
$ws.Range('J435').Value = 'A'
$ws.Range('K435').Value = 2.6
$ws.Range('L435').Value = 3.2
$ws.Range('M435').Value = 2.75
$ws.Range('N435').Value = 2.5
$ws.Range('O435').Value = 3.2
$ws.Range('P435').Value = 3
$ws.Range('Q435').Value = -0.25
$ws.Range('R435').Value = 2.075
$ws.Range('S435').Value = 1.725
$ws.Range('T435').Value = 2
$ws.Range('U435').Value = 1.875
$ws.Range('V435').Value = 1.975
$ws.Range('W435').Value = -1
$ws.Range('X435').Value = -1
$ws.Range('Y435').Value = 2
$ws.Range('Z435').Value = -1
$ws.Range('AA435').Value = 0.7250000000000001
$ws.Range('AB435').Value = 0.875
$ws.Range('AC435').Value = -1

# Row 442
$ws.Range('B442').Value = 5169324
$ws.Range('F442').Value = 'Hibernian'
$ws.Range('G442').Value = 'Livingston'
$ws.Range('H442').Value = 4
$ws.Range('I442').Value = 0
$ws.Range('J442').Value = 'H'
$ws.Range('K442').Value = 2.1
$ws.Range('L442').Value = 3.25
$ws.Range('M442').Value = 3.6
$ws.Range('N442').Value = 1.666
$ws.Range('O442').Value = 4
$ws.Range('P442').Value = 4.8
$ws.Range('Q442').Value = -0.75
$ws.Range('R442').Value = 1.875
$ws.Range('S442').Value = 1.975
$ws.Range('T442').Value = 2.5
$ws.Range('U442').Value = 1.95
$ws.Range('V442').Value = 1.9
$ws.Range('W442').Value = 0.6659999999999999
$ws.Range('X442').Value = -1
$ws.Range('Y442').Value = -1
$ws.Range('Z442').Value = 0.875
$ws.Range('AA442').Value = -1
$ws.Range('AB442').Value = 0.95
$ws.Range('AC442').Value = -1

# Row 443
$ws.Range('B443').Value = 5169864
$ws.Range('F443').Value = 'Dundee Utd'
$ws.Range('G443').Value = 'Hearts'
$ws.Range('H443').Value = 2
$ws.Range('I443').Value = 2
$ws.Range('J443').Value = 'D'
$ws.Range('K443').Value = 3.75
$ws.Range('L443').Value = 3.4
$ws.Range('M443').Value = 2
$ws.Range('N443').Value = 3.1
$ws.Range('O443').Value = 3.6
$ws.Range('P443').Value = 2.15
$ws.Range('Q443').Value = 0.25
$ws.Range('R443').Value = 1.975
$ws.Range('S443').Value = 1.875
$ws.Range('T443').Value = 2.75
$ws.Range('U443').Value = 2
$ws.Range('V443').Value = 1.85
$ws.Range('W443').Value = -1
$ws.Range('X443').Value = 2.6
$ws.Range('Y443').Value = -1
$ws.Range('Z443').Value = 0.4875
$ws.Range('AA443').Value = -0.5
$ws.Range('AB443').Value = 1
$ws.Range('AC443').Value = -1

# Row 444
$ws.Range('B444').Value = 5169865
$ws.Range('F444').Value = 'Rangers'
$ws.Range('G444').Value = 'Motherwell'
$ws.Range('H444').Value = 3
$ws.Range('I444').Value = 0
$ws.Range('J444').Value = 'H'
$ws.Range('K444').Value = 1.25
$ws.Range('L444').Value = 5.5
$ws.Range('M444').Value = 12
$ws.Range('N444').Value = 1.2
$ws.Range('O444').Value = 8
$ws.Range('P444').Value = 10
$ws.Range('Q444').Value = -2
$ws.Range('R444').Value = 2.05
$ws.Range('S444').Value = 1.8
$ws.Range('T444').Value = 3.25
$ws.Range('U444').Value = 1.95
$ws.Range('V444').Value = 1.9
$ws.Range('W444').Value = 0.2
$ws.Range('X444').Value = -1
$ws.Range('Y444').Value = -1
$ws.Range('Z444').Value = 1.05
$ws.Range('AA444').Value = -1
$ws.Range('AB444').Value = -0.5
$ws.Range('AC444').Value = 0.45

# Row 445
$ws.Range('B445').Value = 5169326
$ws.Range('F445').Value = 'Livingston'
$ws.Range('G445').Value = 'St Mirren'
$ws.Range('H445').Value = 1
$ws.Range('I445').Value = 1
$ws.Range('J445').Value = 'D'
$ws.Range('K445').Value = 2.4
$ws.Range('L445').Value = 3.4
$ws.Range('M445').Value = 2.8
$ws.Range('N445').Value = 2.8
$ws.Range('O445').Value = 3.2
$ws.Range('P445').Value = 2.625
$ws.Range('Q445').Value = 0
$ws.Range('R445').Value = 2.025
$ws.Range('S445').Value = 1.825
$ws.Range('T445').Value = 2.25
$ws.Range('U445').Value = 1.975
$ws.Range('V445').Value = 1.875
$ws.Range('W445').Value = -1
$ws.Range('X445').Value = 2.2
$ws.Range('Y445').Value = -1
$ws.Range('Z445').Value = 0
$ws.Range('AA445').Value = -0
$ws.Range('AB445').Value = -0.5
$ws.Range('AC445').Value = 0.4375

# Row 446
$ws.Range('B446').Value = 5169634
$ws.Range('F446').Value = 'Kilmarnock'
$ws.Range('G446').Value = 'Aberdeen'
$ws.Range('H446').Value = 2
$ws.Range('I446').Value = 1
$ws.Range('J446').Value = 'H'
$ws.Range('K446').Value = 3.2
$ws.Range('L446').Value = 3.4
$ws.Range('M446').Value = 2.2
$ws.Range('N446').Value = 3
$ws.Range('O446').Value = 3.6
$ws.Range('P446').Value = 2.25
$ws.Range('Q446').Value = 0.25
$ws.Range('R446').Value = 1.925
$ws.Range('S446').Value = 1.925
$ws.Range('T446').Value = 2.5
$ws.Range('U446').Value = 1.95
$ws.Range('V446').Value = 1.9
$ws.Range('W446').Value = 2
$ws.Range('X446').Value = -1
$ws.Range('Y446').Value = -1
$ws.Range('Z446').Value = 0.925
$ws.Range('AA446').Value = -1
$ws.Range('AB446').Value = 0.95
$ws.Range('AC446').Value = -1

# Row 447
$ws.Range('B447').Value = 5169632
$ws.Range('F447').Value = 'Dundee Utd'
$ws.Range('G447').Value = 'Ross County'
$ws.Range('H447').Value = 3
$ws.Range('I447').Value = 0
$ws.Range('J447').Value = 'H'
$ws.Range('K447').Value = 2.6
$ws.Range('L447').Value = 3.2
$ws.Range('M447').Value = 2.75
$ws.Range('N447').Value = 1.909
$ws.Range('O447').Value = 3.4
$ws.Range('P447').Value = 4.2
$ws.Range('Q447').Value = -0.5
$ws.Range('R447').Value = 1.925
$ws.Range('S447').Value = 1.925
$ws.Range('T447').Value = 2.25
$ws.Range('U447').Value = 1.95
$ws.Range('V447').Value = 1.9
$ws.Range('W447').Value = 0.909
$ws.Range('X447').Value = -1
$ws.Range('Y447').Value = -1
$ws.Range('Z447').Value = 0.925
$ws.Range('AA447').Value = -1
$ws.Range('AB447').Value = 0.95
$ws.Range('AC447').Value = -1

# Row 448
$ws.Range('B448').Value = 5169327
$ws.Range('F448').Value = 'St Johnstone'
$ws.Range('G448').Value = 'Hearts'
$ws.Range('H448').Value = 2
$ws.Range('I448').Value = 3
$ws.Range('J448').Value = 'A'
$ws.Range('K448').Value = 2.7
$ws.Range('L448').Value = 3.25
$ws.Range('M448').Value = 2.6
$ws.Range('N448').Value = 3.1
$ws.Range('O448').Value = 3.6
$ws.Range('P448').Value = 2.2
$ws.Range('Q448').Value = 0.25
$ws.Range('R448').Value = 1.95
$ws.Range('S448').Value = 1.9
$ws.Range('T448').Value = 2.5
$ws.Range('U448').Value = 1.925
$ws.Range('V448').Value = 1.925
$ws.Range('W448').Value = -1
$ws.Range('X448').Value = -1
$ws.Range('Y448').Value = 1.2
$ws.Range('Z448').Value = -1
$ws.Range('AA448').Value = 0.8999999999999999
$ws.Range('AB448').Value = 0.925
$ws.Range('AC448').Value = -1

# Row 456
$ws.Range('B456').Value = 5169332
$ws.Range('F456').Value = 'Ross County'
$ws.Range('G456').Value = 'Livingston'
$ws.Range('H456').Value = 0
$ws.Range('I456').Value = 2
$ws.Range('J456').Value = 'A'
$ws.Range('K456').Value = 2.6
$ws.Range('L456').Value = 3.1
$ws.Range('M456').Value = 2.7
$ws.Range('N456').Value = 2.7
$ws.Range('O456').Value = 3.1
$ws.Range('P456').Value = 2.8
$ws.Range('Q456').Value = 0
$ws.Range('R456').Value = 1.875
$ws.Range('S456').Value = 1.975
$ws.Range('T456').Value = 2
$ws.Range('U456').Value = 1.95
$ws.Range('V456').Value = 1.9
$ws.Range('W456').Value = -1
$ws.Range('X456').Value = -1
$ws.Range('Y456').Value = 1.8
$ws.Range('Z456').Value = -1
$ws.Range('AA456').Value = 0.9750000000000001
$ws.Range('AB456').Value = 0
$ws.Range('AC456').Value = -0

# Row 457
$ws.Range('B457').Value = 5169333
$ws.Range('F457').Value = 'St Mirren'
$ws.Range('G457').Value = 'Hearts'
$ws.Range('H457').Value = 1
$ws.Range('I457').Value = 1
$ws.Range('J457').Value = 'D'
$ws.Range('K457').Value = 3.5
$ws.Range('L457').Value = 3.1
$ws.Range('M457').Value = 2.15
$ws.Range('N457').Value = 2.7
$ws.Range('O457').Value = 3.3
$ws.Range('P457').Value = 2.6
$ws.Range('Q457').Value = 0
$ws.Range('R457').Value = 1.975
$ws.Range('S457').Value = 1.875
$ws.Range('T457').Value = 2.5
$ws.Range('U457').Value = 1.85
$ws.Range('V457').Value = 2
$ws.Range('W457').Value = -1
$ws.Range('X457').Value = 2.3
$ws.Range('Y457').Value = -1
$ws.Range('Z457').Value = 0
$ws.Range('AA457').Value = -0
$ws.Range('AB457').Value = -1
$ws.Range('AC457').Value = 1

# Row 458
$ws.Range('B458').Value = 5169331
$ws.Range('F458').Value = 'Aberdeen'
$ws.Range('G458').Value = 'St Johnstone'
$ws.Range('H458').Value = 2
$ws.Range('I458').Value = 0
$ws.Range('J458').Value = 'H'
$ws.Range('K458').Value = 2.1
$ws.Range('L458').Value = 3.25
$ws.Range('M458').Value = 3.5
$ws.Range('N458').Value = 1.615
$ws.Range('O458').Value = 3.6
$ws.Range('P458').Value = 5.5
$ws.Range('Q458').Value = -0.75
$ws.Range('R458').Value = 1.875
$ws.Range('S458').Value = 1.975
$ws.Range('T458').Value = 2.5
$ws.Range('U458').Value = 2.025
$ws.Range('V458').Value = 1.825
$ws.Range('W458').Value = 0.615
$ws.Range('X458').Value = -1
$ws.Range('Y458').Value = -1
$ws.Range('Z458').Value = 0.875
$ws.Range('AA458').Value = -1
$ws.Range('AB458').Value = -1
$ws.Range('AC458').Value = 0.825

# Row 466
$ws.Range('B466').Value = 5169639
$ws.Range('F466').Value = 'Hearts'
$ws.Range('G466').Value = 'Aberdeen'
$ws.Range('H466').Value = 5
$ws.Range('I466').Value = 0
$ws.Range('J466').Value = 'H'
$ws.Range('K466').Value = 2.375
$ws.Range('L466').Value = 3.25
$ws.Range('M466').Value = 2.9
$ws.Range('N466').Value = 2.1
$ws.Range('O466').Value = 3.5
$ws.Range('P466').Value = 3.5
$ws.Range('Q466').Value = -0.25
$ws.Range('R466').Value = 1.8
$ws.Range('S466').Value = 2.05
$ws.Range('T466').Value = 2.75
$ws.Range('U466').Value = 1.925
$ws.Range('V466').Value = 1.925
$ws.Range('W466').Value = 1.1
$ws.Range('X466').Value = -1
$ws.Range('Y466').Value = -1
$ws.Range('Z466').Value = 0.8
$ws.Range('AA466').Value = -1
$ws.Range('AB466').Value = 0.925
$ws.Range('AC466').Value = -1

# Row 467
$ws.Range('B467').Value = 5169334
$ws.Range('F467').Value = 'Celtic'
$ws.Range('G467').Value = 'St Mirren'
$ws.Range('H467').Value = 4
$ws.Range('I467').Value = 0
$ws.Range('J467').Value = 'H'
$ws.Range('K467').Value = 1.111
$ws.Range('L467').Value = 8
$ws.Range('M467').Value = 21
$ws.Range('N467').Value = 1.1
$ws.Range('O467').Value = 9
$ws.Range('P467').Value = 26
$ws.Range('Q467').Value = -2.25
$ws.Range('R467').Value = 1.85
$ws.Range('S467').Value = 2
$ws.Range('T467').Value = 3.5
$ws.Range('U467').Value = 2
$ws.Range('V467').Value = 1.85
$ws.Range('W467').Value = 0.1000000000000001
$ws.Range('X467').Value = -1
$ws.Range('Y467').Value = -1
$ws.Range('Z467').Value = 0.8500000000000001
$ws.Range('AA467').Value = -1
$ws.Range('AB467').Value = 1
$ws.Range('AC467').Value = -1

# Row 471
$ws.Range('B471').Value = 5169643
$ws.Range('F471').Value = 'Hibernian'
$ws.Range('G471').Value = 'Aberdeen'
$ws.Range('H471').Value = 6
$ws.Range('I471').Value = 0
$ws.Range('J471').Value = 'H'
$ws.Range('K471').Value = 2.1
$ws.Range('L471').Value = 3.4
$ws.Range('M471').Value = 3.4
$ws.Range('N471').Value = 2.25
$ws.Range('O471').Value = 3.6
$ws.Range('P471').Value = 3
$ws.Range('Q471').Value = -0.25
$ws.Range('R471').Value = 2.025
$ws.Range('S471').Value = 1.825
$ws.Range('T471').Value = 2.5
$ws.Range('U471').Value = 1.925
$ws.Range('V471').Value = 1.925
$ws.Range('W471').Value = 1.25
$ws.Range('X471').Value = -1
$ws.Range('Y471').Value = -1
$ws.Range('Z471').Value = 1.025
$ws.Range('AA471').Value = -1
$ws.Range('AB471').Value = 0.925
$ws.Range('AC471').Value = -1

# Row 472
$ws.Range('B472').Value = 5169644
$ws.Range('F472').Value = 'Ross County'
$ws.Range('G472').Value = 'Kilmarnock'
$ws.Range('H472').Value = 3
$ws.Range('I472').Value = 0
$ws.Range('J472').Value = 'H'
$ws.Range('K472').Value = 2.4
$ws.Range('L472').Value = 3.1
$ws.Range('M472').Value = 3.1
$ws.Range('N472').Value = 2.8
$ws.Range('O472').Value = 3
$ws.Range('P472').Value = 2.75
$ws.Range('Q472').Value = 0
$ws.Range('R472').Value = 2
$ws.Range('S472').Value = 1.85
$ws.Range('T472').Value = 2
$ws.Range('U472').Value = 1.9
$ws.Range('V472').Value = 1.95
$ws.Range('W472').Value = 1.8
$ws.Range('X472').Value = -1
$ws.Range('Y472').Value = -1
$ws.Range('Z472').Value = 1
$ws.Range('AA472').Value = -1
$ws.Range('AB472').Value = 0.8999999999999999
$ws.Range('AC472').Value = -1

# Row 478
$ws.Range('B478').Value = 5169340
$ws.Range('F478').Value = 'Celtic'
$ws.Range('G478').Value = 'Livingston'
$ws.Range('H478').Value = 3
$ws.Range('I478').Value = 0
$ws.Range('J478').Value = 'H'
$ws.Range('K478').Value = 1.125
$ws.Range('L478').Value = 8
$ws.Range('M478').Value = 21
$ws.Range('N478').Value = 1.125
$ws.Range('O478').Value = 8.5
$ws.Range('P478').Value = 21
$ws.Range('Q478').Value = -2.5
$ws.Range('R478').Value = 1.975
$ws.Range('S478').Value = 1.875
$ws.Range('T478').Value = 3.75
$ws.Range('U478').Value = 2
$ws.Range('V478').Value = 1.85
$ws.Range('W478').Value = 0.125
$ws.Range('X478').Value = -1
$ws.Range('Y478').Value = -1
$ws.Range('Z478').Value = 0.9750000000000001
$ws.Range('AA478').Value = -1
$ws.Range('AB478').Value = -1
$ws.Range('AC478').Value = 0.8500000000000001

# Row 479
$ws.Range('B479').Value = 5169341
$ws.Range('F479').Value = 'Motherwell'
$ws.Range('G479').Value = 'St Johnstone'
$ws.Range('H479').Value = 0
$ws.Range('I479').Value = 2
$ws.Range('J479').Value = 'A'
$ws.Range('K479').Value = 2.25
$ws.Range('L479').Value = 3.2
$ws.Range('M479').Value = 3.2
$ws.Range('N479').Value = 1.909
$ws.Range('O479').Value = 3.3
$ws.Range('P479').Value = 4.5
$ws.Range('Q479').Value = -0.5
$ws.Range('R479').Value = 1.925
$ws.Range('S479').Value = 1.925
$ws.Range('T479').Value = 2.25
$ws.Range('U479').Value = 1.95
$ws.Range('V479').Value = 1.9
$ws.Range('W479').Value = -1
$ws.Range('X479').Value = -1
$ws.Range('Y479').Value = 3.5
$ws.Range('Z479').Value = -1
$ws.Range('AA479').Value = 0.925
$ws.Range('AB479').Value = -0.5
$ws.Range('AC479').Value = 0.45

# Row 483
$ws.Range('B483').Value = 5169344
$ws.Range('F483').Value = 'St Mirren'
$ws.Range('G483').Value = 'Hibernian'
$ws.Range('H483').Value = 0
$ws.Range('I483').Value = 1
$ws.Range('J483').Value = 'A'
$ws.Range('K483').Value = 2.8
$ws.Range('L483').Value = 3.2
$ws.Range('M483').Value = 2.5
$ws.Range('N483').Value = 2.45
$ws.Range('O483').Value = 3.1
$ws.Range('P483').Value = 2.9
$ws.Range('Q483').Value = -0.25
$ws.Range('R483').Value = 2.075
$ws.Range('S483').Value = 1.725
$ws.Range('T483').Value = 2.25
$ws.Range('U483').Value = 1.95
$ws.Range('V483').Value = 1.9
$ws.Range('W483').Value = -1
$ws.Range('X483').Value = -1
$ws.Range('Y483').Value = 1.9
$ws.Range('Z483').Value = -1
$ws.Range('AA483').Value = 0.7250000000000001
$ws.Range('AB483').Value = -1
$ws.Range('AC483').Value = 0.8999999999999999

# Row 484
$ws.Range('B484').Value = 5169870
$ws.Range('F484').Value = 'Hearts'
$ws.Range('G484').Value = 'Dundee Utd'
$ws.Range('H484').Value = 3
$ws.Range('I484').Value = 1
$ws.Range('J484').Value = 'H'
$ws.Range('K484').Value = 1.533
$ws.Range('L484').Value = 4
$ws.Range('M484').Value = 6
$ws.Range('N484').Value = 1.7
$ws.Range('O484').Value = 3.75
$ws.Range('P484').Value = 5
$ws.Range('Q484').Value = -0.75
$ws.Range('R484').Value = 1.875
$ws.Range('S484').Value = 1.975
$ws.Range('T484').Value = 2.75
$ws.Range('U484').Value = 1.95
$ws.Range('V484').Value = 1.9
$ws.Range('W484').Value = 0.7
$ws.Range('X484').Value = -1
$ws.Range('Y484').Value = -1
$ws.Range('Z484').Value = 0.875
$ws.Range('AA484').Value = -1
$ws.Range('AB484').Value = 0.95
$ws.Range('AC484').Value = -1

# Row 485
$ws.Range('B485').Value = 5169647
$ws.Range('F485').Value = 'Aberdeen'
$ws.Range('G485').Value = 'Motherwell'
$ws.Range('H485').Value = 3
$ws.Range('I485').Value = 1
$ws.Range('J485').Value = 'H'
$ws.Range('K485').Value = 1.909
$ws.Range('L485').Value = 3.4
$ws.Range('M485').Value = 4
$ws.Range('N485').Value = 1.8
$ws.Range('O485').Value = 3.5
$ws.Range('P485').Value = 4.2
$ws.Range('Q485').Value = -0.5
$ws.Range('R485').Value = 1.85
$ws.Range('S485').Value = 2
$ws.Range('T485').Value = 2.5
$ws.Range('U485').Value = 1.95
$ws.Range('V485').Value = 1.9
$ws.Range('W485').Value = 0.8
$ws.Range('X485').Value = -1
$ws.Range('Y485').Value = -1
$ws.Range('Z485').Value = 0.8500000000000001
$ws.Range('AA485').Value = -1
$ws.Range('AB485').Value = 0.95
$ws.Range('AC485').Value = -1

# Row 488
$ws.Range('B488').Value = 5169649
$ws.Range('F488').Value = 'Celtic'
$ws.Range('G488').Value = 'Aberdeen'
$ws.Range('H488').Value = 4
$ws.Range('I488').Value = 0
$ws.Range('J488').Value = 'H'
$ws.Range('K488').Value = 1.166
$ws.Range('L488').Value = 7.5
$ws.Range('M488').Value = 17
$ws.Range('N488').Value = 1.166
$ws.Range('O488').Value = 7.5
$ws.Range('P488').Value = 17
$ws.Range('Q488').Value = -2.25
$ws.Range('R488').Value = 1.975
$ws.Range('S488').Value = 1.875
$ws.Range('T488').Value = 3.5
$ws.Range('U488').Value = 1.9
$ws.Range('V488').Value = 1.95
$ws.Range('W488').Value = 0.1659999999999999
$ws.Range('X488').Value = -1
$ws.Range('Y488').Value = -1
$ws.Range('Z488').Value = 0.9750000000000001
$ws.Range('AA488').Value = -1
$ws.Range('AB488').Value = 0.8999999999999999
$ws.Range('AC488').Value = -1

# Row 489
$ws.Range('B489').Value = 5169346
$ws.Range('F489').Value = 'Livingston'
$ws.Range('G489').Value = 'Rangers'
$ws.Range('H489').Value = 0
$ws.Range('I489').Value = 3
$ws.Range('J489').Value = 'A'
$ws.Range('K489').Value = 8
$ws.Range('L489').Value = 5
$ws.Range('M489').Value = 1.363
$ws.Range('N489').Value = 8
$ws.Range('O489').Value = 5
$ws.Range('P489').Value = 1.363
$ws.Range('Q489').Value = 1.5
$ws.Range('R489').Value = 1.8
$ws.Range('S489').Value = 2.05
$ws.Range('T489').Value = 2.75
$ws.Range('U489').Value = 1.925
$ws.Range('V489').Value = 1.925
$ws.Range('W489').Value = -1
$ws.Range('X489').Value = -1
$ws.Range('Y489').Value = 0.363
$ws.Range('Z489').Value = -1
$ws.Range('AA489').Value = 1.05
$ws.Range('AB489').Value = 0.4625
$ws.Range('AC489').Value = -0.5

# Row 490
$ws.Range('B490').Value = 5169345
$ws.Range('F490').Value = 'Dundee Utd'
$ws.Range('G490').Value = 'St Johnstone'
$ws.Range('H490').Value = 1
$ws.Range('I490').Value = 2
$ws.Range('J490').Value = 'A'
$ws.Range('K490').Value = 2.25
$ws.Range('L490').Value = 3.25
$ws.Range('M490').Value = 3.25
$ws.Range('N490').Value = 2.45
$ws.Range('O490').Value = 3.2
$ws.Range('P490').Value = 3
$ws.Range('Q490').Value = -0.25
$ws.Range('R490').Value = 2.1
$ws.Range('S490').Value = 1.775
$ws.Range('T490').Value = 2.25
$ws.Range('U490').Value = 2.05
$ws.Range('V490').Value = 1.75
$ws.Range('W490').Value = -1
$ws.Range('X490').Value = -1
$ws.Range('Y490').Value = 2
$ws.Range('Z490').Value = -1
$ws.Range('AA490').Value = 0.7749999999999999
$ws.Range('AB490').Value = 1.05
$ws.Range('AC490').Value = -1

# Row 533
$ws.Range('B533').Value = 5428458
$ws.Range('F533').Value = 'St Johnstone'
$ws.Range('G533').Value = 'Hibernian'
$ws.Range('H533').Value = 1
$ws.Range('I533').Value = 1
$ws.Range('J533').Value = 'D'
$ws.Range('K533').Value = 3.1
$ws.Range('L533').Value = 3.25
$ws.Range('M533').Value = 2.3
$ws.Range('N533').Value = 3.8
$ws.Range('O533').Value = 3.3
$ws.Range('P533').Value = 2.05
$ws.Range('Q533').Value = 0.5
$ws.Range('R533').Value = 1.8
$ws.Range('S533').Value = 2.05
$ws.Range('T533').Value = 2.25
$ws.Range('U533').Value = 1.875
$ws.Range('V533').Value = 1.975
$ws.Range('W533').Value = -1
$ws.Range('X533').Value = 2.3
$ws.Range('Y533').Value = -1
$ws.Range('Z533').Value = 0.8
$ws.Range('AA533').Value = -1
$ws.Range('AB533').Value = -0.5
$ws.Range('AC533').Value = 0.4875

# Row 534
$ws.Range('B534').Value = 5498114
$ws.Range('F534').Value = 'Dundee Utd'
$ws.Range('G534').Value = 'Livingston'
$ws.Range('H534').Value = 2
$ws.Range('I534').Value = 0
$ws.Range('J534').Value = 'H'
$ws.Range('K534').Value = 2.5
$ws.Range('L534').Value = 3.2
$ws.Range('M534').Value = 2.8
$ws.Range('N534').Value = 2.7
$ws.Range('O534').Value = 3.3
$ws.Range('P534').Value = 2.75
$ws.Range('Q534').Value = 0
$ws.Range('R534').Value = 1.9
$ws.Range('S534').Value = 1.95
$ws.Range('T534').Value = 2.25
$ws.Range('U534').Value = 1.9
$ws.Range('V534').Value = 1.95
$ws.Range('W534').Value = 1.7
$ws.Range('X534').Value = -1
$ws.Range('Y534').Value = -1
$ws.Range('Z534').Value = 0.8999999999999999
$ws.Range('AA534').Value = -1
$ws.Range('AB534').Value = -0.5
$ws.Range('AC534').Value = 0.475

# Row 605
$ws.Range('B605').Value = 6844740
$ws.Range('F605').Value = 'St Johnstone'
$ws.Range('G605').Value = 'Livingston'
$ws.Range('H605').Value = 1
$ws.Range('I605').Value = 1
$ws.Range('J605').Value = 'D'
$ws.Range('K605').Value = 2.5
$ws.Range('L605').Value = 3.3
$ws.Range('M605').Value = 2.8
$ws.Range('N605').Value = 2.4
$ws.Range('O605').Value = 3.1
$ws.Range('P605').Value = 3.1
$ws.Range('Q605').Value = -0.25
$ws.Range('R605').Value = 2.05
$ws.Range('S605').Value = 1.8
$ws.Range('T605').Value = 2
$ws.Range('U605').Value = 1.875
$ws.Range('V605').Value = 1.975
$ws.Range('W605').Value = -1
$ws.Range('X605').Value = 2.1
$ws.Range('Y605').Value = -1
$ws.Range('Z605').Value = -0.5
$ws.Range('AA605').Value = 0.4
$ws.Range('AB605').Value = 0
$ws.Range('AC605').Value = -0

# Row 606
$ws.Range('B606').Value = 6845506
$ws.Range('F606').Value = 'Rangers'
$ws.Range('G606').Value = 'Aberdeen'
$ws.Range('H606').Value = 1
$ws.Range('I606').Value = 3
$ws.Range('J606').Value = 'A'
$ws.Range('K606').Value = 1.333
$ws.Range('L606').Value = 5
$ws.Range('M606').Value = 9
$ws.Range('N606').Value = 1.363
$ws.Range('O606').Value = 4.75
$ws.Range('P606').Value = 8
$ws.Range('Q606').Value = -1.5
$ws.Range('R606').Value = 2
$ws.Range('S606').Value = 1.85
$ws.Range('T606').Value = 3
$ws.Range('U606').Value = 1.975
$ws.Range('V606').Value = 1.875
$ws.Range('W606').Value = -1
$ws.Range('X606').Value = -1
$ws.Range('Y606').Value = 7
$ws.Range('Z606').Value = -1
$ws.Range('AA606').Value = 0.8500000000000001
$ws.Range('AB606').Value = 0.9750000000000001
$ws.Range('AC606').Value = -1

# Row 654
$ws.Range('B654').Value = 6844781
$ws.Range('F654').Value = 'St Johnstone'
$ws.Range('G654').Value = 'St Mirren'
$ws.Range('H654').Value = 1
$ws.Range('I654').Value = 0
$ws.Range('J654').Value = 'H'
$ws.Range('K654').Value = 2.7
$ws.Range('L654').Value = 3.2
$ws.Range('M654').Value = 2.6
$ws.Range('N654').Value = 3
$ws.Range('O654').Value = 2.9
$ws.Range('P654').Value = 2.625
$ws.Range('Q654').Value = 0
$ws.Range('R654').Value = 1.925
$ws.Range('S654').Value = 1.925
$ws.Range('T654').Value = 2
$ws.Range('U654').Value = 2.05
$ws.Range('V654').Value = 1.8
$ws.Range('W654').Value = 2
$ws.Range('X654').Value = -1
$ws.Range('Y654').Value = -1
$ws.Range('Z654').Value = 0.925
$ws.Range('AA654').Value = -1
$ws.Range('AB654').Value = -1
$ws.Range('AC654').Value = 0.8

# Row 655
$ws.Range('B655').Value = 6844778
$ws.Range('F655').Value = 'Aberdeen'
$ws.Range('G655').Value = 'Kilmarnock'
$ws.Range('H655').Value = 0
$ws.Range('I655').Value = 1
$ws.Range('J655').Value = 'A'
$ws.Range('K655').Value = 2.25
$ws.Range('L655').Value = 3.25
$ws.Range('M655').Value = 3.2
$ws.Range('N655').Value = 2.1
$ws.Range('O655').Value = 3.4
$ws.Range('P655').Value = 3.6
$ws.Range('Q655').Value = -0.5
$ws.Range('R655').Value = 2.05
$ws.Range('S655').Value = 1.8
$ws.Range('T655').Value = 2.25
$ws.Range('U655').Value = 1.9
$ws.Range('V655').Value = 1.95
$ws.Range('W655').Value = -1
$ws.Range('X655').Value = -1
$ws.Range('Y655').Value = 2.6
$ws.Range('Z655').Value = -1
$ws.Range('AA655').Value = 0.8
$ws.Range('AB655').Value = -1
$ws.Range('AC655').Value = 0.95
